$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) stores numeric-looking values (e.g. "215.18") as plain
# text. Assigning such a string straight to .Value lets Excel silently
# coerce it to a real number, so each touched Price cell is flipped to the
# Text number format just long enough to take the literal string, then
# restored to the workbooks normal (General) style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.797.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.626.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5115"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.247"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.634.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.849.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7494"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.806.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.427"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.775"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.013"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.880"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1249"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.714"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.239"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04874"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.170"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.360"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8957"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5525"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.116.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01549"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.548"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7960"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.775.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4424"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9989"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.613"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.16%  "
